$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared string used across the whole "source" column (F2:F29)
# from "houndstooth" to "Elvis" for every cell that carries that string.
$ws.Range("F2:F29").Value = "Elvis"

# Updated avg_sent_length (D) / lexical_diversity (E) values for the rows whose
# underlying text files changed after quote-stripping.
$ws.Range("D2").Value = 129.9990862364363
$ws.Range("E2").Value = 0.2115677056328275
$ws.Range("D9").Value = 84.6236483309826
$ws.Range("E9").Value = 0.2091244451498848
$ws.Range("D10").Value = 75.86591412596266
$ws.Range("E10").Value = 0.2133620483369227
$ws.Range("D11").Value = 70.99365009853295
$ws.Range("E11").Value = 0.2115136010656817
$ws.Range("D12").Value = 84.07059306382105
$ws.Range("E12").Value = 0.2116515206824339
$ws.Range("D13").Value = 82.18699380657456
$ws.Range("E13").Value = 0.2074510870180924
$ws.Range("D14").Value = 88.84892086330935
$ws.Range("E14").Value = 0.2064019923820685
$ws.Range("D15").Value = 82.27863569222505
$ws.Range("E15").Value = 0.2102729866591266
$ws.Range("D18").Value = 88.789631405454
$ws.Range("E18").Value = 0.2080328255874656
$ws.Range("D21").Value = 118.9483695652174
$ws.Range("E21").Value = 0.2118452172598859
$ws.Range("D22").Value = 126.2081542571727
$ws.Range("E22").Value = 0.2006435728703712
$ws.Range("D23").Value = 129.2213103352344
$ws.Range("E23").Value = 0.2019204377945379
$ws.Range("D25").Value = 136.5766764858777
$ws.Range("E25").Value = 0.2077986601830865
$ws.Range("D26").Value = 103.8485148514851
$ws.Range("E26").Value = 0.212412310601038
$ws.Range("D28").Value = 112.9669443697931
$ws.Range("E28").Value = 0.2106708227257448
$ws.Range("D29").Value = 105.6900746474979
$ws.Range("E29").Value = 0.2107714791805455
